# Applies the "old word review" feature data change:
# Rows 6 through 14 (the mekong river delta .. the red river delta) have their
# Level (column C) changed from "nothing" to "easy".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 6; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = "easy"
}
